$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.099.85"
$ws.Range("E2").Value = "  +11.55%  "
$ws.Range("D3").Value = "1.813.83"
$ws.Range("E3").Value = "  +8.06%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'228.17"
$ws.Range("E5").Value = "  +3.66%  "
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'31.40"
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("D9").Value = "'46.66"
$ws.Range("E9").Value = "  +5.24%  "
$ws.Range("D10").Value = "'0.281"
$ws.Range("E10").Value = "  +6.13%  "
$ws.Range("E11").Value = "  +4.96%  "
$ws.Range("D12").Value = "'0.0927"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "2.074.38"
$ws.Range("E13").Value = "  +8.03%  "
$ws.Range("D14").Value = "1.810.10"
$ws.Range("E14").Value = "  +7.72%  "
$ws.Range("D15").Value = "'0.640"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").Value = "34.038.30"
$ws.Range("E16").Value = "  +11.31%  "
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "'4.27"
$ws.Range("E18").Value = "  +7.24%  "
$ws.Range("D19").Value = "'69.35"
$ws.Range("E19").Value = "  +4.43%  "
$ws.Range("D20").Value = "'257.59"
$ws.Range("E20").Value = "  +5.10%  "
$ws.Range("D21").Value = "0.0₃0746"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").Value = "'158.27"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'16.58"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("E28").Value = "  +6.26%  "
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.87"
$ws.Range("E31").Value = "  +10.86%  "
$ws.Range("B32").Value = "MinaProtocolToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D32").Value = "'1.79"
$ws.Range("E32").Value = "  +334.03%  "
$ws.Range("D33").Value = "'0.0512"
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("D35").Value = "'3.52"
$ws.Range("E35").Value = "  +6.57%  "
$ws.Range("D36").Value = "1.540.80"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("D41").Value = "'0.622"
$ws.Range("E41").Value = "  +5.05%  "
$ws.Range("D42").Value = "'2.80"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "'0.908"
$ws.Range("E44").Value = "  +8.29%  "
$ws.Range("E45").Value = "  +8.00%  "
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("E47").Value = "  +4.40%  "
$ws.Range("D48").Value = "1.969.02"
$ws.Range("E48").Value = "  +8.35%  "
$ws.Range("D49").Value = "'5.71"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "'52.59"
$ws.Range("E51").Value = "  +1.71%  "
